# Saldo.xlsx update:
#  - two balances increase slightly
#  - a handful of new account rows are inserted, in the right spots,
#    so the sheet stays sorted by descending balance
#  - two new (negative) HFR rows are appended at the bottom of the data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update two existing balances (rows are untouched by any insert,
#        so their row numbers don't move) ---
$ws.Range("C2").Value = 63757.2    # 008012870 / ANA
$ws.Range("C3").Value = 63756.64   # 008004995 / JOSE

# --- 2) Insert the new rows, from top to bottom, using each row's FINAL
#        (post-insert) row number. Because we go top to bottom, every
#        later insertion's target row number already accounts for the
#        rows inserted above it. ---

function Add-Row($rowNum, $conta, $nome, $saldo) {
    $ws.Rows.Item($rowNum).Insert()
    $ws.Range("A$rowNum").NumberFormat = "@"
    $ws.Range("A$rowNum").Value = $conta
    $ws.Range("B$rowNum").Value = $nome
    $ws.Range("C$rowNum").Value = $saldo
}

Add-Row 4   "005547703" "SILVIA"   35168
Add-Row 5   "004468717" "HELOISA"  23079
Add-Row 7   "003512801" "LAIS"     3115.52
Add-Row 20  "004231371" "ADRIANO"  500
Add-Row 30  "004487016" "ROGERIO"  184.88
Add-Row 92  "005514036" "ANA"      41.34

# --- 3) Now fix MIRELLA's balance. Her row was originally row 7, but by
#        this point three rows have been inserted above it (4, 5, 7), so
#        she now sits at row 10. ---
$ws.Range("C10").Value = 931.65

# --- 4) Append the two new HFR rows right after the last data row
#        (NORTON, account 004976625) and before the existing trailing
#        blank row. All six earlier inserts happened above NORTON's
#        original row 180, so NORTON is now at row 186, and the new HFR
#        rows land at 187 and 188. ---
Add-Row 187 "004381180" "HFR" -27335.61
Add-Row 188 "004361159" "HFR" -300454.07
